$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.321.11"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "2.356.11"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "540.23"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").Value = "135.45"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  +5.10%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "5.58"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "0.353"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.775.21"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "23.80"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "58.318.59"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "2.357.04"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "10.72"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "332.81"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "63.04"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "8.47"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "1.39"
$ws.Range("E27").Value = "  +3.70%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "172.02"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.75"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "1.04"
$ws.Range("E32").Value = "  +13.55%  "
$ws.Range("D33").Value = "18.45"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +6.73%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("D39").Value = "39.24"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "145.62"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "294.52"
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("D42").Value = "0.379"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").Value = "0.0948"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").Value = "19.19"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").Value = "0.386"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "17.55"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +0.62%  "
